$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.770.93"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.745.69"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.89"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.75"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.111"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.390"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.53"
$ws.Range("E12").Value = "  -18.63%  "
$ws.Range("D13").Value = "3.229.45"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.95"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "63.675.00"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("D17").Value = "2.753.50"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.27"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "360.06"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.564"
$ws.Range("E22").Value = "  +4.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.990"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.32"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "0.0₃0933"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  -2.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.11"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.35"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  +4.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.49"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.12"
$ws.Range("E40").Value = "  +7.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "332.02"
$ws.Range("E41").Value = "  -5.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.62"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.96"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0258"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.638"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.67"
$ws.Range("E48").Value = "  -4.74%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("E51").Value = "  +1.10%  "
